$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add "yes" values to C2, D2, E2 (style matches surrounding row via style index 1 / the 14pt font style already applied on the sheet)
$ws.Range("C2").Value = "yes"
$ws.Range("D2").Value = "yes"
$ws.Range("E2").Value = "yes"

# Move the active selection from B2 to C7
$ws.Range("C7").Select()
